$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a daily price list (newest date at row 2, oldest at the bottom).
# A new day's entry (17-11-2025) is being added at the top, pushing every
# existing row down by one. The row that used to be the last one (row 159,
# 12-06-2025) becomes a brand-new row 160.

$lastOldRow = 159
$lastNewRow = 160

# Shift all data rows down by one, working from the bottom up so we never
# overwrite a source row before it has been copied.
for ($n = $lastNewRow; $n -ge 3; $n--) {
    $src = $n - 1
    $ws.Range("A" + $src + ":F" + $src).Copy($ws.Range("A" + $n + ":F" + $n))
}

# Row 2 keeps the same description/grade/price/circular info as before (the
# price did not change between 16-11-2025 and 17-11-2025); only the date
# changes. Copy an existing text cell into A2 first so the new value keeps
# its text (inline/shared-string) type instead of being auto-converted to a
# date serial number, then overwrite with the new date.
$ws.Range("A3").Copy($ws.Range("A2"))
$ws.Range("A2").Value2 = "17-11-2025"

# Rebuild every hyperlink on column F from scratch, since copying cells does
# not bring hyperlinks along with it in this engine.
$ws.Hyperlinks.Delete()
for ($n = 2; $n -le $lastNewRow; $n++) {
    $link = $ws.Range("F" + $n).Value2
    if ($link -ne $null -and $link -ne "") {
        $ws.Hyperlinks.Add($ws.Range("F" + $n), $link) | Out-Null
    }
}
